$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.915.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4600"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3878"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07866"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9859"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.897.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.994"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.644"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06961"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009981"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.911.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.242"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.098"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.054"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.927"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09341"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9032"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.317"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.255"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.180"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05752"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.669"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5650"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.666"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.274"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5345"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07045"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.841"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.528"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.059"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.26%  "
$ws.Range("E51").Value = "  -0.04%  "
